$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bsg"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 38.60286033333333
$ws.Range("H2").Value2 = 115.808581
$ws.Range("I2").Value2 = 0.2650212684862838
$ws.Range("J2").Value2 = 0.2650212684862838
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 7.879565666666667
$ws.Range("N2").Value2 = 23.638697
$ws.Range("O2").Value2 = 0.9977172793687663
$ws.Range("P2").Value2 = 0.9977172793687664
$ws.Range("Q2").Value2 = 304.1737729176619
$ws.Range("R2").Value2 = 2737.563956258957
$ws.Range("S2").Value2 = 0.2644162989689944
$ws.Range("T2").Value2 = 0.2644162989689945

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bsg"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 38.60286033333333
$ws.Range("H3").Value2 = 115.808581
$ws.Range("I3").Value2 = 0.2650212684862838
$ws.Range("J3").Value2 = 0.2650212684862838
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 0.3333333333333333
$ws.Range("M3").Value2 = 0.018028
$ws.Range("N3").Value2 = 0.054084
$ws.Range("O3").Value2 = 0.002282720631233623
$ws.Range("P3").Value2 = 0.002282720631233623
$ws.Range("Q3").Value2 = 0.6959323660893333
$ws.Range("R3").Value2 = 6.263391294804
$ws.Range("S3").Value2 = 0.0006049695172893454
$ws.Range("T3").Value2 = 0.0006049695172893454

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bsg"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 74.45592499999999
$ws.Range("H4").Value2 = 223.367775
$ws.Range("I4").Value2 = 0.5111642898850374
$ws.Range("J4").Value2 = 0.5111642898850374
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 7.879565666666667
$ws.Range("N4").Value2 = 23.638697
$ws.Range("O4").Value2 = 0.9977172793687663
$ws.Range("P4").Value2 = 0.9977172793687664
$ws.Range("Q4").Value2 = 586.6803503099084
$ws.Range("R4").Value2 = 5280.123152789175
$ws.Range("S4").Value2 = 0.5099974446145669
$ws.Range("T4").Value2 = 0.509997444614567

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bsg"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 74.45592499999999
$ws.Range("H5").Value2 = 223.367775
$ws.Range("I5").Value2 = 0.5111642898850374
$ws.Range("J5").Value2 = 0.5111642898850374
$ws.Range("K5").Value2 = 1
$ws.Range("L5").Value2 = 0.3333333333333333
$ws.Range("M5").Value2 = 0.018028
$ws.Range("N5").Value2 = 0.054084
$ws.Range("O5").Value2 = 0.002282720631233623
$ws.Range("P5").Value2 = 0.002282720631233623
$ws.Range("Q5").Value2 = 1.3422914159
$ws.Range("R5").Value2 = 12.0806227431
$ws.Range("S5").Value2 = 0.001166845270470459
$ws.Range("T5").Value2 = 0.001166845270470459

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Bsg"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 32.60069533333333
$ws.Range("H6").Value2 = 97.80208599999999
$ws.Range("I6").Value2 = 0.2238144416286788
$ws.Range("J6").Value2 = 0.2238144416286788
$ws.Range("K6").Value2 = 3
$ws.Range("L6").Value2 = 1
$ws.Range("M6").Value2 = 7.879565666666667
$ws.Range("N6").Value2 = 23.638697
$ws.Range("O6").Value2 = 0.9977172793687663
$ws.Range("P6").Value2 = 0.9977172793687664
$ws.Range("Q6").Value2 = 256.8793196579935
$ws.Range("R6").Value2 = 2311.913876921942
$ws.Range("S6").Value2 = 0.2233035357852049
$ws.Range("T6").Value2 = 0.223303535785205

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Bsg"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 32.60069533333333
$ws.Range("H7").Value2 = 97.80208599999999
$ws.Range("I7").Value2 = 0.2238144416286788
$ws.Range("J7").Value2 = 0.2238144416286788
$ws.Range("K7").Value2 = 1
$ws.Range("L7").Value2 = 0.3333333333333333
$ws.Range("M7").Value2 = 0.018028
$ws.Range("N7").Value2 = 0.054084
$ws.Range("O7").Value2 = 0.002282720631233623
$ws.Range("P7").Value2 = 0.002282720631233623
$ws.Range("Q7").Value2 = 0.5877253354693331
$ws.Range("R7").Value2 = 5.289528019223999
$ws.Range("S7").Value2 = 0.0005109058434738185
$ws.Range("T7").Value2 = 0.0005109058434738185
